# Apply the edit described by the commit: correct the "Nápojové kartony"
# (beverage cartons) weight figure for the "Paskov - zástavba rodinných domů"
# village (cell N25) from 2.9 to -2.9, then leave the selection on the merged
# cell P4:Q4 as the last user action (matching the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Core data edit: N25 changes sign (2.9 -> -2.9). All dependent formulas
# (O/T/U columns, row 28/29/31/32/34 totals, etc.) recalculate automatically.
$ws.Range("N25").Value = -2.9

# Leave the active selection on the merged date-header cell for the
# "Rudolfov" column (P4:Q4), matching the final saved workbook/view state
# (this also clears the previous topLeftCell="A10" scroll position).
$ws.Range("P4:Q4").Select()
